# Applies the "Updated cryptos list" GitHub Actions data refresh to Sheet1.
# Rows 2-51 hold a scraped Coin/Link/Price/Volume(1h) snapshot; this run's
# diff updates Price/Volume cells throughout, plus rows 24-27 where the
# scrape order shifted so Coin/Link/Price/Volume moved to a different row.
#
# Price/Volume text (columns D/E) looks numeric ("1.896.18", "  +0.49%  ")
# but is stored as plain text in the source data, not real numbers/percents.
# A leading apostrophe forces Excel to keep it as text instead of silently
# parsing it into a Number (which would mangle the double-dot price format
# and drop significant trailing zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.217.95"
$ws.Range("E2").Value = "'  +0.49%  "
$ws.Range("D3").Value = "'1.896.50"
$ws.Range("E3").Value = "'  -0.03%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "'  +0.21%  "
$ws.Range("D5").Value = "'307.47"
$ws.Range("E5").Value = "'  +0.09%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "'  +0.13%  "
$ws.Range("D7").Value = "'0.5202"
$ws.Range("E7").Value = "'  -0.59%  "
$ws.Range("D8").Value = "'0.3765"
$ws.Range("E8").Value = "'  -0.28%  "
$ws.Range("D9").Value = "'0.07286"
$ws.Range("E9").Value = "'  +0.87%  "
$ws.Range("E10").Value = "'  +0.11%  "
$ws.Range("D11").Value = "'0.9004"
$ws.Range("E11").Value = "'  +1.01%  "
$ws.Range("D12").Value = "'0.08170"
$ws.Range("E12").Value = "'  +6.40%  "
$ws.Range("D13").Value = "'96.63"
$ws.Range("E13").Value = "'  +2.44%  "
$ws.Range("D14").Value = "'1.893.94"
$ws.Range("E14").Value = "'  -0.12%  "
$ws.Range("D15").Value = "'5.283"
$ws.Range("E15").Value = "'  +1.02%  "
$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "'  +0.12%  "
$ws.Range("D17").Value = "'0.000008613"
$ws.Range("E17").Value = "'  +1.21%  "
$ws.Range("D18").Value = "'14.55"
$ws.Range("E18").Value = "'  +0.34%  "
$ws.Range("D20").Value = "'27.249.93"
$ws.Range("D21").Value = "'5.087"
$ws.Range("E21").Value = "'  +0.44%  "
$ws.Range("E22").Value = "'  +0.97%  "
$ws.Range("D23").Value = "'6.407"
$ws.Range("B24").Value = "Monero"
$ws.Range("C24").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D24").Value = "'147.66"
$ws.Range("E24").Value = "'  +1.32%  "
$ws.Range("B25").Value = "LidoDAOToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D25").Value = "'2.294"
$ws.Range("E25").Value = "'  +0.19%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "'1.748"
$ws.Range("E26").Value = "'  +0.76%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'18.22"
$ws.Range("E27").Value = "'  +0.78%  "
$ws.Range("D28").Value = "'115.09"
$ws.Range("E28").Value = "'  +0.45%  "
$ws.Range("D29").Value = "'4.957"
$ws.Range("E29").Value = "'  -0.05%  "
$ws.Range("D30").Value = "'4.834"
$ws.Range("E30").Value = "'  +0.84%  "
$ws.Range("D31").Value = "'0.09223"
$ws.Range("E31").Value = "'  +0.14%  "
$ws.Range("D32").Value = "'0.05033"
$ws.Range("E32").Value = "'  -0.29%  "
$ws.Range("D33").Value = "'0.7946"
$ws.Range("E33").Value = "'  +2.28%  "
$ws.Range("D34").Value = "'1.219"
$ws.Range("E34").Value = "'  -1.60%  "
$ws.Range("D35").Value = "'3.448"
$ws.Range("E35").Value = "'  +4.55%  "
$ws.Range("D36").Value = "'2.944"
$ws.Range("E36").Value = "'  -1.08%  "
$ws.Range("D37").Value = "'2.591"
$ws.Range("E37").Value = "'  -0.01%  "
$ws.Range("D38").Value = "'0.5671"
$ws.Range("E38").Value = "'  +0.42%  "
$ws.Range("E39").Value = "'  -0.24%  "
$ws.Range("E40").Value = "'  +0.17%  "
$ws.Range("E41").Value = "'  -0.27%  "
$ws.Range("D42").Value = "'6.559"
$ws.Range("E42").Value = "'  -1.05%  "
$ws.Range("D43").Value = "'115.43"
$ws.Range("E43").Value = "'  -2.65%  "
$ws.Range("E44").Value = "'  -0.15%  "
$ws.Range("D45").Value = "'0.4898"
$ws.Range("E45").Value = "'  +1.33%  "
$ws.Range("E46").Value = "'  +0.10%  "
$ws.Range("E47").Value = "'  -0.65%  "
$ws.Range("D48").Value = "'1.623"
$ws.Range("E48").Value = "'  +1.67%  "
$ws.Range("D49").Value = "'38.21"
$ws.Range("E49").Value = "'  +1.93%  "
$ws.Range("D50").Value = "'63.45"
$ws.Range("E50").Value = "'  -0.90%  "
$ws.Range("D51").Value = "'0.05939"
$ws.Range("E51").Value = "'  +0.32%  "
